# Added Login Scenarios Dated 26052024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "LoggedInUser" locator row (row 22) to the ObjectRepository sheet.
# Shared strings must be introduced in this column order (D, B, A) so the new
# unique strings land at shared-string indices 47, 48, 49 respectively,
# matching the target workbook.
$ws.Range("D22").Value = "//a[@class='account']//span"
$ws.Range("B22").Value = "LoggedInUser"
$ws.Range("C22").Value = "xpath"
$ws.Range("A22").Value = "HomePage"

# Match styling of the other data rows (centered horizontal/vertical alignment,
# i.e. the same cell style already used throughout column A:D).
$ws.Range("A22:D22").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A22:D22").VerticalAlignment = -4108    # xlCenter

# Update the active selection left behind by the editing session.
$null = $ws.Range("F12").Select()
